$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PDF Prefix Registry")

# Insert a new row at position 5 (shifts existing rows 5-75 down to 6-76,
# carrying along cell styles/row heights as Excel natively does).
$ws.Rows.Item(5).Insert()

# Populate the new row 5 with the Affinitext Inc. / Andras Dippold entry
# (Issue #34: https://github.com/adobe/pdf-names-list/issues/34)
$ws.Range("A5").Value = "aftx"
$ws.Range("B5").Value = "Andras"
$ws.Range("C5").Value = "Dippold"
$ws.Range("D5").Value = "andras.dippold@affinitext.com"
$ws.Range("E5").Value = "Affinitext Inc."
$ws.Range("F5").Value = 44987
$ws.Range("G5").Value = "https://github.com/adobe/pdf-names-list/issues/34 "

# The row-insert does not re-point the worksheet's existing Hyperlinks
# collection at their new (shifted) rows, so rebuild it from scratch:
# drop every old hyperlink, then re-add one per data row at its new
# location (mailto targets taken from the pre-edit relationships).
$ws.Range("A1").Hyperlinks.Delete()

$hyperlinkRows = @(
    @(2, "mailto:frank.vyncke@yin4yang.com"),
    @(3, "mailto:infohelp@websupergoo.com"),
    @(4, "mailto:doberndorfer@adlibsoftware.com"),
    @(6, "mailto:steve@ajidev.com"),
    @(7, "mailto:geojohn@att.net"),
    @(8, "mailto:joerg@bec.de"),
    @(9, "mailto:mike@bfo.com"),
    @(10, "mailto:wagner@bipro.net"),
    @(11, "mailto:callas@compuserve.com"),
    @(12, "mailto:qianxiaobo@cfca.com.cn"),
    @(13, "mailto:qianxiaobo@cfca.com.cn"),
    @(15, "mailto:info@stefanochizzolini.it"),
    @(16, "mailto:mstorer@autonomy.com"),
    @(17, "mailto:m.schubert@disphere.com"),
    @(18, "mailto:boris.doubrov@duallab.com"),
    @(19, "mailto:horvath.istvan@dmsone.hu"),
    @(20, "mailto:duff.johnson@pdfa.org"),
    @(21, "mailto:lieven.plettinck@esko.com"),
    @(22, "mailto:rob.evans@domuslocus.com"),
    @(23, "mailto:bertvk@enfocus.com"),
    @(24, "mailto:raffael.boesch@power.alstom.com"),
    @(25, "mailto:sbingham@fileopen.com"),
    @(26, "mailto:roman_toda@foxitsoftware.com"),
    @(27, "mailto:dstrebe@mapthematics.com"),
    @(28, "mailto:andreas@goofyfootsoftware.com"),
    @(29, "mailto:martin.bailey@globalgraphics.com"),
    @(30, "mailto:q@as207960.net"),
    @(31, "mailto:betsy.fanning@3dpdfconsortium.com"),
    @(32, "mailto:dorf@aptech.com"),
    @(33, "mailto:alexei@gurnetgroup.com"),
    @(34, "mailto:david@zwang.com"),
    @(35, "mailto:smith.michael@hebcom"),
    @(36, "mailto:frankv@hybridsoftware.com"),
    @(37, "mailto:ymjeong@icerti.com"),
    @(38, "mailto:s.chantaraud@interactivepdf.org"),
    @(39, "mailto:betsy.fanning@3dpdfconsortium.com"),
    @(40, "mailto:bruno@1t3xt.com"),
    @(41, "mailto:ueharak@kthree.co.jp"),
    @(42, "mailto:adobe@kawseq.com"),
    @(43, "mailto:d.luksaite@archyvai.lt"),
    @(44, "mailto:klink@csi.com"),
    @(45, "mailto:barakc@microsoft.com"),
    @(46, "mailto:abdias.lira@wolterskluwer.com"),
    @(47, "mailto:mek@mekentosj.com"),
    @(48, "mailto:toda@digitaldocuments.org"),
    @(49, "mailto:tvtreeck@nepatec.de"),
    @(50, "mailto:tvtreeck@nepatec.de"),
    @(51, "mailto:opena@its.jnj.com"),
    @(52, "mailto:wilson.hirata@iti.gov.br"),
    @(53, "mailto:duff.johnson@pdfa.org"),
    @(54, "mailto:pdfthingys@shaw.ca"),
    @(55, "mailto:pierre.choutet@laposte.net"),
    @(56, "mailto:tm@pdflib.com"),
    @(57, "mailto:pragma@wxs.nl"),
    @(58, "mailto:info@sense4code.com"),
    @(59, "mailto:brett.mckenzie@scinaptic.com"),
    @(60, "mailto:jan.slabon@setasign.com"),
    @(61, "mailto:info@sejda.org"),
    @(62, "mailto:normand.williams@onespan.com"),
    @(63, "mailto:gwarren@snl.com"),
    @(64, "mailto:coeters@sofha.de"),
    @(65, "mailto:stephen@stephenreindl.net"),
    @(66, "mailto:betsy.fanning@3dpdfconsortium.com"),
    @(67, "mailto:hcobb@telegenisys.com"),
    @(68, "mailto:daniel.banks@citycomputers.co.uk"),
    @(69, "mailto:John_Brinkema@ao.uscourts.gov"),
    @(70, "mailto:silviu.ardelean@visma.com"),
    @(71, "mailto:nick.mettyear@wacom.eu"),
    @(72, "mailto:nick.mettyear@wacom.eu"),
    @(73, "mailto:colin.borrowman@wolterskluwer.com"),
    @(74, "mailto:vaidas@superita.it"),
    @(75, "mailto:nathan@xmpie.com"),
    @(76, "mailto:grzegorz.molenda@zeto.lublin.pl")

)

foreach ($pair in $hyperlinkRows) {
    $r = $pair[0]
    $addr = $pair[1]
    $ws.Hyperlinks.Add($ws.Range("D" + $r), $addr) | Out-Null
    $ws.Range("D" + $r).Style = "Hyperlink"
}

# New hyperlinks introduced by this edit
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:andras.dippold@affinitext.com") | Out-Null
$ws.Range("D5").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("G5"), "https://github.com/adobe/pdf-names-list/issues/34 ") | Out-Null
$ws.Range("G5").Style = "Hyperlink"
